# Update Name of Algo
# Applies updated RandomForest-imputed values to the BC/20/seed5 result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = -12.3283
$ws.Range("C4").Value = -12.0677
$ws.Range("B11").Value = 5.792799999999999
$ws.Range("B12").Value = 5.135299999999997
$ws.Range("C14").Value = -12.9523
$ws.Range("B15").Value = 4.498799999999996
$ws.Range("C26").Value = -11.3747
$ws.Range("B27").Value = 6.42
$ws.Range("B28").Value = 6.547900000000002
$ws.Range("B31").Value = 3.841699999999995
$ws.Range("C31").Value = -13.4593
$ws.Range("B32").Value = 6.661399999999999
$ws.Range("C35").Value = -12.08830000000001
$ws.Range("B36").Value = 9.121499999999999
$ws.Range("C37").Value = -12.9605
$ws.Range("B38").Value = 4.966099999999998
$ws.Range("C39").Value = -12.79960000000001
$ws.Range("C40").Value = -13.39800000000001
$ws.Range("C45").Value = -13.2677
$ws.Range("B46").Value = 7.446800000000003
$ws.Range("C52").Value = -11.0462
$ws.Range("B54").Value = 4.8326
$ws.Range("B55").Value = 4.834499999999998
$ws.Range("B56").Value = 4.495299999999998
$ws.Range("C57").Value = -14.68639999999999
$ws.Range("B67").Value = 5.304299999999997
$ws.Range("B69").Value = 5.604299999999994
$ws.Range("B72").Value = 5.836899999999998
$ws.Range("B73").Value = 8.416699999999995
$ws.Range("C81").Value = -12.88079999999999
$ws.Range("B83").Value = 5.081799999999994
$ws.Range("C83").Value = -14.0228
$ws.Range("B86").Value = 4.738200000000004
$ws.Range("B91").Value = 5.2747
$ws.Range("B93").Value = 6.544300000000004
$ws.Range("B99").Value = 4.787799999999995
$ws.Range("C100").Value = -12.9072
$ws.Range("C102").Value = -12.939
